# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) on several sheets to freshly
# pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2507
$ws.Range("J17").Value = 2507
$ws.Range("L17").Value = 7521
$ws.Range("N17").Value = -7857

$ws.Range("H28").Value = 1098.6666
$ws.Range("I28").Value = 416.33334
$ws.Range("K28").Value = 416.33334
$ws.Range("M28").Value = 68.66665999999998

$ws.Range("H70").Value = 3133.2222
$ws.Range("J70").Value = 4166.5
$ws.Range("L70").Value = 12499.5
$ws.Range("N70").Value = -13039.5

$ws.Range("H73").Value = 3133.2222
$ws.Range("J73").Value = 4166.5
$ws.Range("L73").Value = 12499.5
$ws.Range("N73").Value = -14371.5

$ws.Range("H76").Value = 3664.8333
$ws.Range("I76").Value = 2663
$ws.Range("K76").Value = 2663
$ws.Range("M76").Value = -2348

$ws.Range("H79").Value = 3664.8333
$ws.Range("I79").Value = 2663
$ws.Range("K79").Value = 2663
$ws.Range("M79").Value = -1571

$ws.Range("H98").Value = 757.75
$ws.Range("I98").Value = 631.06665
$ws.Range("K98").Value = 631.06665
$ws.Range("M98").Value = 866.93335

$ws.Range("H116").Value = 4473.75
$ws.Range("I116").Value = 3965
$ws.Range("K116").Value = 3965
$ws.Range("M116").Value = -523

$ws.Range("H122").Value = 757.75
$ws.Range("I122").Value = 631.06665
$ws.Range("K122").Value = 1893.19995
$ws.Range("M122").Value = 556.8000500000001

$ws.Range("H125").Value = 3172.7856
$ws.Range("I125").Value = 2035
$ws.Range("J125").Value = 9999.5
$ws.Range("K125").Value = 18315
$ws.Range("L125").Value = 89995.5
$ws.Range("M125").Value = -15855
$ws.Range("N125").Value = -94915.5

$ws.Range("H138").Value = 3366.4211
$ws.Range("I138").Value = 2018.3334
$ws.Range("J138").Value = 4579.7
$ws.Range("K138").Value = 6055.0002
$ws.Range("L138").Value = 13739.1
$ws.Range("M138").Value = -915.0002000000004
$ws.Range("N138").Value = -24019.1

$ws.Range("H141").Value = 3106.8125
$ws.Range("I141").Value = 3106.8125
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9320.4375
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4140.4375
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 8753
$ws.Range("I26").Value = 8753
$ws.Range("K26").Value = 8753
$ws.Range("M26").Value = -8423

$ws.Range("H39").Value = 14999.5
$ws.Range("I39").Value = 9999
$ws.Range("K39").Value = 9999
$ws.Range("M39").Value = -9479

$ws.Range("H61").Value = 3400
$ws.Range("I61").Value = 3400
$ws.Range("K61").Value = 3400
$ws.Range("M61").Value = -3188

$ws.Range("H74").Value = 4649.5
$ws.Range("I74").Value = 3976.111
$ws.Range("J74").Value = 6669.6665
$ws.Range("K74").Value = 3976.111
$ws.Range("L74").Value = 6669.6665
$ws.Range("M74").Value = -3102.111
$ws.Range("N74").Value = -8417.666499999999

$ws.Range("H77").Value = 4649.5
$ws.Range("I77").Value = 3976.111
$ws.Range("J77").Value = 6669.6665
$ws.Range("K77").Value = 19880.555
$ws.Range("L77").Value = 33348.3325
$ws.Range("M77").Value = -15512.555
$ws.Range("N77").Value = -42084.3325

$ws.Range("H92").Value = 41633.332
$ws.Range("J92").Value = 41633.332
$ws.Range("L92").Value = 41633.332
$ws.Range("N92").Value = -46625.332

$ws.Range("H132").Value = 1807.875
$ws.Range("I132").Value = 1807.875
$ws.Range("K132").Value = 5423.625
$ws.Range("M132").Value = -2893.625

$ws.Range("H136").Value = 3400
$ws.Range("I136").Value = 3400
$ws.Range("K136").Value = 10200
$ws.Range("M136").Value = -7650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4599.6665
$ws.Range("I134").Value = 4599.6665
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13798.9995
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11263.9995
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H17").Value = 13987.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 13987.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 13987.25
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -14335.25

$ws.Range("H19").Value = 578.8182
$ws.Range("I19").Value = 144.375
$ws.Range("J19").Value = 1737.3334
$ws.Range("K19").Value = 144.375
$ws.Range("L19").Value = 1737.3334
$ws.Range("M19").Value = 25.625
$ws.Range("N19").Value = -2077.3334

$ws.Range("H24").Value = 578.8182
$ws.Range("I24").Value = 144.375
$ws.Range("J24").Value = 1737.3334
$ws.Range("K24").Value = 144.375
$ws.Range("L24").Value = 1737.3334
$ws.Range("M24").Value = 25.625
$ws.Range("N24").Value = -2077.3334

$ws.Range("H86").Value = 12621.6
$ws.Range("I86").Value = 12406.667
$ws.Range("K86").Value = 12406.667
$ws.Range("M86").Value = -11283.667

$ws.Range("H89").Value = 12621.6
$ws.Range("I89").Value = 12406.667
$ws.Range("K89").Value = 62033.335
$ws.Range("M89").Value = -56417.335

$ws.Range("H112").Value = 59999.5
$ws.Range("J112").Value = 59999.5
$ws.Range("L112").Value = 59999.5
$ws.Range("N112").Value = -62953.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8700

$ws.Range("H102").Value = 4987.4
$ws.Range("I102").Value = 4987.4
$ws.Range("K102").Value = 4987.4
$ws.Range("M102").Value = -3365.4

$ws.Range("H122").Value = 3593.5557
$ws.Range("I122").Value = 2269
$ws.Range("K122").Value = 6807
$ws.Range("M122").Value = -4357

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4999.5
$ws.Range("I40").Value = 4999.5
$ws.Range("K40").Value = 4999.5
$ws.Range("M40").Value = -4863.5

$ws.Range("H46").Value = 3559.8
$ws.Range("I46").Value = 2600
$ws.Range("J46").Value = 4076.6155
$ws.Range("K46").Value = 2600
$ws.Range("L46").Value = 4076.6155
$ws.Range("M46").Value = -2412
$ws.Range("N46").Value = -4452.6155

$ws.Range("H122").Value = 3252.5
$ws.Range("I122").Value = 3252.5
$ws.Range("K122").Value = 9757.5
$ws.Range("M122").Value = -7307.5

$ws.Range("H136").Value = 3502.3333
$ws.Range("I136").Value = 3502.3333
$ws.Range("K136").Value = 10506.9999
$ws.Range("M136").Value = -7956.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H113").Value = 9101.083000000001
$ws.Range("I113").Value = 17340
$ws.Range("K113").Value = 52020
$ws.Range("M113").Value = -49850

$ws.Range("H122").Value = 3851
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800

$ws.Range("H126").Value = 882.5
$ws.Range("I126").Value = 979.3
$ws.Range("K126").Value = 2937.9
$ws.Range("M126").Value = -467.8999999999996

$ws.Range("H132").Value = 1771.3572
$ws.Range("I132").Value = 1771.3572
$ws.Range("K132").Value = 5314.071599999999
$ws.Range("M132").Value = -2784.071599999999
